# Auto-generated edit script: apply numeric updates to Leve profit tables
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 34000
$ws.Cells.Item(124, 10).Value = 34000
$ws.Cells.Item(124, 12).Value = 34000
$ws.Cells.Item(124, 14).Value = -43820
$ws.Cells.Item(126, 8).Value = 14979
$ws.Cells.Item(126, 10).Value = 14979
$ws.Cells.Item(126, 12).Value = 14979
$ws.Cells.Item(126, 14).Value = -24859
$ws.Cells.Item(135, 8).Value = 1681.52
$ws.Cells.Item(135, 9).Value = 1650.6957
$ws.Cells.Item(135, 10).Value = 2036
$ws.Cells.Item(135, 11).Value = 14856.2613
$ws.Cells.Item(135, 12).Value = 18324
$ws.Cells.Item(135, 13).Value = -12321.2613
$ws.Cells.Item(135, 14).Value = -23394

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3887.0625
$ws.Cells.Item(61, 9).Value = 3097.75
$ws.Cells.Item(61, 10).Value = 6255
$ws.Cells.Item(61, 11).Value = 3097.75
$ws.Cells.Item(61, 12).Value = 6255
$ws.Cells.Item(61, 13).Value = -2885.75
$ws.Cells.Item(61, 14).Value = -6679
$ws.Cells.Item(74, 8).Value = 9708.532999999999
$ws.Cells.Item(74, 9).Value = 1200
$ws.Cells.Item(74, 10).Value = 17153.5
$ws.Cells.Item(74, 11).Value = 1200
$ws.Cells.Item(74, 12).Value = 17153.5
$ws.Cells.Item(74, 13).Value = -326
$ws.Cells.Item(74, 14).Value = -18901.5
$ws.Cells.Item(77, 8).Value = 9708.532999999999
$ws.Cells.Item(77, 9).Value = 1200
$ws.Cells.Item(77, 10).Value = 17153.5
$ws.Cells.Item(77, 11).Value = 6000
$ws.Cells.Item(77, 12).Value = 85767.5
$ws.Cells.Item(77, 13).Value = -1632
$ws.Cells.Item(77, 14).Value = -94503.5
$ws.Cells.Item(122, 8).Value = 6453
$ws.Cells.Item(122, 9).Value = 5648
$ws.Cells.Item(122, 10).Value = 7580
$ws.Cells.Item(122, 11).Value = 16944
$ws.Cells.Item(122, 12).Value = 22740
$ws.Cells.Item(122, 13).Value = -14494
$ws.Cells.Item(122, 14).Value = -27640
$ws.Cells.Item(132, 8).Value = 3850.3684
$ws.Cells.Item(132, 9).Value = 3409.8
$ws.Cells.Item(132, 10).Value = 5502.5
$ws.Cells.Item(132, 11).Value = 10229.4
$ws.Cells.Item(132, 12).Value = 16507.5
$ws.Cells.Item(132, 13).Value = -7699.400000000001
$ws.Cells.Item(132, 14).Value = -21567.5
$ws.Cells.Item(136, 8).Value = 3887.0625
$ws.Cells.Item(136, 9).Value = 3097.75
$ws.Cells.Item(136, 10).Value = 6255
$ws.Cells.Item(136, 11).Value = 9293.25
$ws.Cells.Item(136, 12).Value = 18765
$ws.Cells.Item(136, 13).Value = -6743.25
$ws.Cells.Item(136, 14).Value = -23865

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3421.5945
$ws.Cells.Item(134, 9).Value = 2493.9167
$ws.Cells.Item(134, 10).Value = 5134.231
$ws.Cells.Item(134, 11).Value = 7481.750100000001
$ws.Cells.Item(134, 12).Value = 15402.693
$ws.Cells.Item(134, 13).Value = -4946.750100000001
$ws.Cells.Item(134, 14).Value = -20472.693

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5933.7715
$ws.Cells.Item(31, 9).Value = 2397.6155
$ws.Cells.Item(31, 10).Value = 8023.3184
$ws.Cells.Item(31, 11).Value = 2397.6155
$ws.Cells.Item(31, 12).Value = 8023.3184
$ws.Cells.Item(31, 13).Value = -2102.6155
$ws.Cells.Item(31, 14).Value = -8613.3184
$ws.Cells.Item(34, 8).Value = 5933.7715
$ws.Cells.Item(34, 9).Value = 2397.6155
$ws.Cells.Item(34, 10).Value = 8023.3184
$ws.Cells.Item(34, 11).Value = 2397.6155
$ws.Cells.Item(34, 12).Value = 8023.3184
$ws.Cells.Item(34, 13).Value = -2195.6155
$ws.Cells.Item(34, 14).Value = -8427.3184
$ws.Cells.Item(134, 8).Value = 3254.2258
$ws.Cells.Item(134, 9).Value = 1746.625
$ws.Cells.Item(134, 10).Value = 4862.3335
$ws.Cells.Item(134, 11).Value = 5239.875
$ws.Cells.Item(134, 12).Value = 14587.0005
$ws.Cells.Item(134, 13).Value = -2704.875
$ws.Cells.Item(134, 14).Value = -19657.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 360.3243
$ws.Cells.Item(107, 10).Value = 366.06668
$ws.Cells.Item(107, 12).Value = 1098.20004
$ws.Cells.Item(107, 14).Value = -4938.20004
$ws.Cells.Item(113, 8).Value = 945.1429000000001
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 945.1429000000001
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 2835.4287
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -7175.4287
$ws.Cells.Item(136, 8).Value = 3747.3
$ws.Cells.Item(136, 9).Value = 725
$ws.Cells.Item(136, 10).Value = 4212.269
$ws.Cells.Item(136, 11).Value = 2175
$ws.Cells.Item(136, 12).Value = 12636.807
$ws.Cells.Item(136, 13).Value = 2925
$ws.Cells.Item(136, 14).Value = -22836.807

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 7409637
$ws.Cells.Item(11, 9).Value = 9055557
$ws.Cells.Item(11, 10).Value = 3000
$ws.Cells.Item(11, 11).Value = 9055557
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = -9055418
$ws.Cells.Item(11, 14).Value = -3278
$ws.Cells.Item(122, 8).Value = 2150.875
$ws.Cells.Item(122, 9).Value = 1857
$ws.Cells.Item(122, 10).Value = 4208
$ws.Cells.Item(122, 11).Value = 5571
$ws.Cells.Item(122, 12).Value = 12624
$ws.Cells.Item(122, 13).Value = -3121
$ws.Cells.Item(122, 14).Value = -17524
$ws.Cells.Item(123, 8).Value = 15934.19
$ws.Cells.Item(123, 10).Value = 15934.19
$ws.Cells.Item(123, 12).Value = 15934.19
$ws.Cells.Item(123, 14).Value = -20834.19
$ws.Cells.Item(126, 8).Value = 2720.5757
$ws.Cells.Item(126, 9).Value = 2341.3572
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 7024.071599999999
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -4554.071599999999
$ws.Cells.Item(126, 14).Value = -13940
$ws.Cells.Item(132, 8).Value = 3367.5715
$ws.Cells.Item(132, 9).Value = 3663.1538
$ws.Cells.Item(132, 10).Value = 2887.25
$ws.Cells.Item(132, 11).Value = 10989.4614
$ws.Cells.Item(132, 12).Value = 8661.75
$ws.Cells.Item(132, 13).Value = -8459.4614
$ws.Cells.Item(132, 14).Value = -13721.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(92, 8).Value = 32000
$ws.Cells.Item(92, 10).Value = 32000
$ws.Cells.Item(92, 12).Value = 32000
$ws.Cells.Item(92, 14).Value = -36992
$ws.Cells.Item(136, 8).Value = 11480.4
$ws.Cells.Item(136, 9).Value = 7134
$ws.Cells.Item(136, 10).Value = 18000
$ws.Cells.Item(136, 11).Value = 21402
$ws.Cells.Item(136, 12).Value = 54000
$ws.Cells.Item(136, 13).Value = -18852
$ws.Cells.Item(136, 14).Value = -59100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 28696.55
$ws.Cells.Item(125, 10).Value = 28696.55
$ws.Cells.Item(125, 12).Value = 28696.55
$ws.Cells.Item(125, 14).Value = -38536.55
$ws.Cells.Item(132, 8).Value = 1689.9791
$ws.Cells.Item(132, 9).Value = 1553.7028
$ws.Cells.Item(132, 10).Value = 2148.3635
$ws.Cells.Item(132, 11).Value = 4661.1084
$ws.Cells.Item(132, 12).Value = 6445.0905
$ws.Cells.Item(132, 13).Value = -2131.1084
$ws.Cells.Item(132, 14).Value = -11505.0905
$ws.Cells.Item(136, 8).Value = 2113.5715
$ws.Cells.Item(136, 9).Value = 1029.6154
$ws.Cells.Item(136, 10).Value = 3875
$ws.Cells.Item(136, 11).Value = 3088.8462
$ws.Cells.Item(136, 12).Value = 11625
$ws.Cells.Item(136, 13).Value = -538.8462
$ws.Cells.Item(136, 14).Value = -16725

